$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teams Templates")

# ------------------------------------------------------------------
# 1. Re-populate the "Teams Templates" table (rows 7-19).
#    The underlying SharePoint query table was refreshed with the
#    newly released Teams templates and the whole list re-sorted in
#    ascending (alphabetical) order, so rows 7-10 get new content and
#    nine brand-new template rows (11-19) are appended.
#    Row 6 carries the canonical formatting (style pattern
#    B=3, C=1, D=1, E=5, F=2, G=2) used by every data row in the
#    table, so it is used as the format donor for every new/changed
#    row before the literal values are written.
# ------------------------------------------------------------------

# Row 7: Organize a store
$ws.Range("A6:G6").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)
$ws.Range("A7").Value = "Organize a store"
$ws.Range("B7").Value = "retailStore"
$ws.Range("C7").Value = "retailStore"
$ws.Range("E7").Value = "Channels:`n`nShift handoff`nLearning`n`nTeam properties`n`nTeam visibility set to Public`n`nMember permissions`n`nPrevent members from creating, updating, or removing channels`nPrevent members from adding or removing apps`nPrevent members from creating, updating, or removing connectors"
$ws.Range("F7").Value = $true
$ws.Rows.Item(7).RowHeight = 246.5

# Row 8: Retail - Manager collaboration
$ws.Range("A6:G6").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Retail - Manager collaboration"
$ws.Range("B8").Value = "retailManagerCollaboration"
$ws.Range("C8").Value = "retailManagerCollaboration"
$ws.Range("E8").Value = "Channels:`n`nShift handoff`nLearning`n`nTeam properties:`n`nTeam visibility set to Private`n`nMember permissions:`n`nPrevent members from creating, updating, or removing channels`nPrevent members from adding or removing apps`nPrevent members from creating, updating, or removing connectors"
$ws.Range("F8").Value = $true
$ws.Rows.Item(8).RowHeight = 246.5

# Row 9: Healthcare - Hospital
$ws.Range("A6:G6").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)
$ws.Range("A9").Value = "Healthcare - Hospital"
$ws.Range("B9").Value = "healthcareHospital"
$ws.Range("C9").Value = "healthcareHospital"
$ws.Range("E9").Value = "Channels:`n`nAnnouncements*`nCompliance*`nCustodial`nHuman Resources`nPharmacy`n`n*Auto-favorited channel"
$ws.Range("F9").Value = $true
$ws.Rows.Item(9).RowHeight = 130.5

# Row 10: Healthcare - Ward
$ws.Range("A6:G6").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Healthcare - Ward"
$ws.Range("B10").Value = "healthcareWard"
$ws.Range("C10").Value = "healthcareWard"
$ws.Range("E10").Value = "Channels:`n`nAnnouncements*`nHuddles*`nRounds`nStaffing*`nTraining*`n`n*Auto-favorited channels"
$ws.Range("F10").Value = $true
$ws.Rows.Item(10).RowHeight = 130.5

# Row 11: Adopt Office 365
$ws.Range("A6:G6").Copy()
$ws.Range("A11:G11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Adopt Office 365"
$ws.Range("B11").Value = "com.microsoft.teams.template.AdoptOffice365"
$ws.Range("C11").Value = "com.microsoft.teams.template.AdoptOffice365"
$ws.Range("E11").Value = "Channels:`nGeneral`nAnnouncements`nChampions corner`nTeam forms`nApps:`nWiki`nCalendar"
$ws.Range("F11").Value = $true
$ws.Rows.Item(11).RowHeight = 159.5

# Row 12: Manage an event
$ws.Range("A6:G6").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Manage an event"
$ws.Range("B12").Value = "com.microsoft.teams.template.ManageAnEvent"
$ws.Range("C12").Value = "com.microsoft.teams.template.ManageAnEvent"
$ws.Range("E12").Value = "Channels:`nGeneral`nAnnouncements`nBudget`nContent`nLogistics`nPlanning`nMarketing and PR`nApps:`nWiki`nWebsite`nYouTube`nPlanner`nOneNote"
$ws.Range("F12").Value = $true
$ws.Rows.Item(12).RowHeight = 246.5

# Row 13: Manage a Project
$ws.Range("A6:G6").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Manage a Project"
$ws.Range("B13").Value = "com.microsoft.teams.template.ManageAProject"
$ws.Range("C13").Value = "com.microsoft.teams.template.ManageAProject"
$ws.Range("E13").Value = "Channels:`nGeneral`nAnnouncements`nResources`nPlanning`nApps:`nWiki`nOneNote"
$ws.Range("F13").Value = $true
$ws.Rows.Item(13).RowHeight = 159.5

# Row 14: Onboard employees
$ws.Range("A6:G6").Copy()
$ws.Range("A14:G14").PasteSpecial(-4122)
$ws.Range("A14").Value = "Onboard employees"
$ws.Range("B14").Value = "com.microsoft.teams.template.OnboardEmployees"
$ws.Range("C14").Value = "com.microsoft.teams.template.OnboardEmployees"
$ws.Range("E14").Value = "Channels:`nGeneral`nAnnouncements`nEmployee chat`nTraining`nApps:`nWiki`nCommunities"
$ws.Range("F14").Value = $true
$ws.Rows.Item(14).RowHeight = 159.5

# Row 15: Organize help desk
$ws.Range("A6:G6").Copy()
$ws.Range("A15:G15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Organize help desk"
$ws.Range("B15").Value = "com.microsoft.teams.template.OrganizeHelpDesk"
$ws.Range("C15").Value = "com.microsoft.teams.template.OrganizeHelpDesk"
$ws.Range("E15").Value = "Channels:`nGeneral`nAnnouncements`nFAQ`nApps:`nWiki`nOneNote"
$ws.Range("F15").Value = $true
$ws.Rows.Item(15).RowHeight = 145

# Row 16: Collaborate on global crisis or event
$ws.Range("A6:G6").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Collaborate on global crisis or event"
$ws.Range("B16").Value = "com.microsoft.teams.template.CollaborateOnAGlobalCrisisOrEvent"
$ws.Range("C16").Value = "com.microsoft.teams.template.CollaborateOnAGlobalCrisisOrEvent"
$ws.Range("E16").Value = "Channels:`nGeneral`nAnnouncements`nWorld news`nBusiness continuity`nRemote working`nInternal comms`nExternal comms`nCustomer complaints`nKudos`nExecutive update`nApps:`nPraise`nWiki`nWebsite"
$ws.Range("F16").Value = $true
$ws.Rows.Item(16).RowHeight = 261

# Row 17: Collaborate within a bank branch
$ws.Range("A6:G6").Copy()
$ws.Range("A17:G17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Collaborate within a bank branch"
$ws.Range("B17").Value = "com.microsoft.teams.template.CollaborateWithinABankBranch"
$ws.Range("C17").Value = "com.microsoft.teams.template.CollaborateWithinABankBranch"
$ws.Range("E17").Value = "Channels:`nGeneral`nAnnouncements`nHuddles`nCustomer meetings`nCoaching`nSkills development`nLoan processing`nCustomer complaints`nKudos`nFun stuff`nCompliance"
$ws.Range("F17").Value = $true
$ws.Rows.Item(17).RowHeight = 188.5

# Row 18: Coordinate incident response
$ws.Range("A6:G6").Copy()
$ws.Range("A18:G18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Coordinate incident response"
$ws.Range("B18").Value = "com.microsoft.teams.template.CoordinateIncidentResponse"
$ws.Range("C18").Value = "com.microsoft.teams.template.CoordinateIncidentResponse"
$ws.Range("E18").Value = "Channels:`nGeneral`nAnnouncements`nLogistics`nPlanning`nRecovery`nUrgent`nApps:`nWiki`nExcel`nOneNote`nSharePoint`nPlanner"
$ws.Range("F18").Value = $true
$ws.Rows.Item(18).RowHeight = 232

# Row 19: Quality and safety
$ws.Range("A6:G6").Copy()
$ws.Range("A19:G19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Quality and safety"
$ws.Range("B19").Value = "com.microsoft.teams.template.QualitySafety"
$ws.Range("C19").Value = "com.microsoft.teams.template.QualitySafety"
$ws.Range("E19").Value = "Channels:`nGeneral`nAnnouncements`nLine 1`nLine 2`nLine 3`nSafety`nTraining`nMaintenance`nFun stuff`nApps:`nWiki"
$ws.Range("F19").Value = $true
$ws.Rows.Item(19).RowHeight = 217.5

$ws.Application.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Grow the query table / ListObject so it covers the new extent.
# ------------------------------------------------------------------
$ws.ListObjects.Item(1).Resize($ws.Range("A1:G19"))

# ------------------------------------------------------------------
# 3. The hidden ExternalData_1 name (source range for the query
#    table) needs to track the new extent too.
# ------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Teams Templates!ExternalData_1") {
        $n.RefersTo = "='Teams Templates'!`$A`$1:`$G`$19"
    }
}

# ------------------------------------------------------------------
# 4. Column A widened to fit the longer template names, and no
#    longer auto "best fit".
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 33.88

# ------------------------------------------------------------------
# 5. Selection moved while reviewing the newly inserted rows.
# ------------------------------------------------------------------
$ws.Range("G12").Select()
